$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 20 was a blank placeholder row in the journal table; fill it in with a
# new entry (date, start/end time, duration, project, description).
$ws.Range("B20").Value = 44258
$ws.Range("C20").Value = "10H40"
$ws.Range("D20").Value = "12h15"
$ws.Range("E20").Value = "95min"
$ws.Range("G20").Value = "Ajout dans la page d'acceuile `n-CSS réorganisé `n-Mis a jour de la version de Bootstrap + `nRéparation des conflits"

# The description text is long, so the row grows taller (matches the other
# multi-line description rows in the sheet).
$ws.Rows.Item(20).RowHeight = 57

# Move the saved scroll position / selection to reflect where the author
# was working.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()
